$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Name "
$ws.Range("B1").Value = "Status"

$ws.Range("B1").Select()
